$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values for the 380 kV case (rows 2-25, data rows 0-23).
# Columns G, K, L, O stay 0 (unchanged) and are intentionally skipped.
$newValues = @{
    2 = @{ "B"=18.67841498217127; "C"=15.58052727591503; "D"=6.381431055820447; "E"=11.55663407063; "F"=44.0447152764466; "H"=7.344005520526261; "I"=28.41591799105309; "J"=9.812307771001171; "M"=19.86841817526788; "N"=19.74095577190074 }
    3 = @{ "B"=18.25681080096359; "C"=15.17686764404522; "D"=6.380924773044715; "E"=11.57567309220718; "F"=43.91941992684603; "H"=7.344005520526261; "I"=28.43656366824172; "J"=9.835241316783421; "M"=19.75657196039053; "N"=19.8129151788179 }
    4 = @{ "B"=17.99843605282652; "C"=14.92853205735604; "D"=6.380833982313594; "E"=11.58867990716314; "F"=43.85473432854631; "H"=7.344005520526261; "I"=28.45626921228926; "J"=9.850360458558169; "M"=19.69225295019293; "N"=19.85908292750368 }
    5 = @{ "B"=17.89343537626833; "C"=14.82737912150118; "D"=6.380852546114456; "E"=11.59431169799893; "F"=43.83146349030735; "H"=7.344005520526261; "I"=28.46606021128455; "J"=9.856782919271541; "M"=19.66715794034377; "N"=19.87839697488432 }
    6 = @{ "B"=17.87602246339823; "C"=14.81059055635248; "D"=6.380858990999593; "E"=11.59526687719415; "F"=43.82778616641863; "H"=7.344005520526261; "I"=28.46779213812418; "J"=9.857865154341702; "M"=19.66305886946195; "N"=19.88163431458894 }
    7 = @{ "B"=17.99701858504258; "C"=14.92716746110031; "D"=6.380834007388833; "E"=11.58875451722243; "F"=43.8544079726228; "H"=7.344005520526261; "I"=28.45639413718213; "J"=9.850446015758365; "M"=19.69190996784502; "N"=19.85934137605486 }
    8 = @{ "B"=18.53304174297647; "C"=15.44154547841037; "D"=6.381210946943813; "E"=11.56292560674217; "F"=43.99897862159582; "H"=7.344005520526261; "I"=28.42157396885901; "J"=9.820000038012115; "M"=19.82896247084893; "N"=19.76535631902192 }
    9 = @{ "B"=19.58105183694634; "C"=16.43933353184716; "D"=6.383687001542272; "E"=11.52271041742152; "F"=44.37906107582295; "H"=7.344005520526261; "I"=28.40932922777305; "J"=9.768516665221979; "M"=20.13130332436365; "N"=19.59673510657795 }
    10 = @{ "B"=20.3401327206374; "C"=17.15691700166241; "D"=6.386553371959587; "E"=11.49950875631161; "F"=44.71616186142798; "H"=7.344005520526261; "I"=28.43480893340833; "J"=9.73568493171828; "M"=20.37249784705768; "N"=19.482321190294 }
    11 = @{ "B"=20.68134180755722; "C"=17.47831832097467; "D"=6.388082347605854; "E"=11.49032751390446; "F"=44.88178729673984; "H"=7.344005520526261; "I"=28.45393521285767; "J"=9.721829280286684; "M"=20.48603514792901; "N"=19.43230925413054 }
    12 = @{ "B"=20.80982493380789; "C"=17.5991748550555; "D"=6.38869347917799; "E"=11.48704795270112; "F"=44.9462404131806; "H"=7.344005520526261; "I"=28.46226373863281; "J"=9.716737480519393; "M"=20.52954822682887; "N"=19.4136623894403 }
    13 = @{ "B"=20.78218811574161; "C"=17.57318601664471; "D"=6.3885604350468; "E"=11.48774550043745; "F"=44.93228267346127; "H"=7.344005520526261; "I"=28.46042172434607; "J"=9.717827198932957; "M"=20.52015431093047; "N"=19.4176653755862 }
    14 = @{ "B"=20.69192747305243; "C"=17.48827901391654; "D"=6.388131982856442; "E"=11.49005375252584; "F"=44.88705532346356; "H"=7.344005520526261; "I"=28.45459863496396; "J"=9.721407268969076; "M"=20.48960476863258; "N"=19.43076932829204 }
    15 = @{ "B"=20.63654189212732; "C"=17.43615656340598; "D"=6.387873722804503; "E"=11.4914932924299; "F"=44.85957715446062; "H"=7.344005520526261; "I"=28.45117327975932; "J"=9.723620348901362; "M"=20.47095897425866; "N"=19.43883381533782 }
    16 = @{ "B"=20.31773958030684; "C"=17.13580027681948; "D"=6.386457959528159; "E"=11.50013638094079; "F"=44.70558231144046; "H"=7.344005520526261; "I"=28.43371083143704; "J"=9.736612137844501; "M"=20.36515218776511; "N"=19.48563045036959 }
    17 = @{ "B"=20.12101496726594; "C"=16.95015867783351; "D"=6.385646918140438; "E"=11.50579016170809; "F"=44.61423479003589; "H"=7.344005520526261; "I"=28.42493008877728; "J"=9.744858549876142; "M"=20.30119975860068; "N"=19.51485920285066 }
    18 = @{ "B"=20.00748574415546; "C"=16.84291597377018; "D"=6.385201624387564; "E"=11.50917134845441; "F"=44.56285171452845; "H"=7.344005520526261; "I"=28.42058883696765; "J"=9.74970330431937; "M"=20.26477711321572; "N"=19.53186246851656 }
    19 = @{ "B"=19.9689859894123; "C"=16.80652936634096; "D"=6.385054502706594; "E"=11.51033837343513; "F"=44.54565399765216; "H"=7.344005520526261; "I"=28.4192406851731; "J"=9.751361117872953; "M"=20.25250794485107; "N"=19.53765243486249 }
    20 = @{ "B"=20.14199681631564; "C"=16.96996978915956; "D"=6.385731062794137; "E"=11.50517492909892; "F"=44.62383929629154; "H"=7.344005520526261; "I"=28.42579139514843; "J"=9.743970187880109; "M"=20.30797044469919; "N"=19.51172792444143 }
    21 = @{ "B"=20.7184599329834; "C"=17.51324232538429; "D"=6.388256959083935; "E"=11.48937041413479; "F"=44.90029288570871; "H"=7.344005520526261; "I"=28.45627953571507; "J"=9.720351509744429; "M"=20.49856406776005; "N"=19.42691247555996 }
    22 = @{ "B"=21.09092764479087; "C"=17.86328917262335; "D"=6.390095023411554; "E"=11.48019039074779; "F"=45.09106381564484; "H"=7.344005520526261; "I"=28.48253462275968; "J"=9.705818901391766; "M"=20.62613848218892; "N"=19.37317965758832 }
    23 = @{ "B"=20.89256905257528; "C"=17.67696078624678; "D"=6.38909695335488; "E"=11.48498490001596; "F"=44.98833338792344; "H"=7.344005520526261; "I"=28.46794217101672; "J"=9.713492623053746; "M"=20.55778425801676; "N"=19.4017027870655 }
    24 = @{ "B"=20.13251225014261; "C"=16.96101477964791; "D"=6.385692955616819; "E"=11.50545266847273; "F"=44.6194935650697; "H"=7.344005520526261; "I"=28.42539979683877; "J"=9.744371493063026; "M"=20.30490834235035; "N"=19.51314295465327 }
    25 = @{ "B"=19.29882785535332; "C"=16.17152155471361; "D"=6.3828324239007; "E"=11.53247424787538; "F"=44.26599034030316; "H"=7.344005520526261; "I"=28.40660658895768; "J"=9.781566112807704; "M"=20.04605660915345; "N"=19.64068186629116 }
}

foreach ($row in $newValues.Keys) {
    foreach ($col in $newValues[$row].Keys) {
        $ws.Range("$col$row").Value = $newValues[$row][$col]
    }
}
